# cmip6 mappings for cmip5-atmos-key-properties.xlsx
# Reproduces the "Added cmip6 mappings" commit: populates the previously
# empty cmip6-id column (B) for the solar / volcanoes rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cmip6-id values in column B -----------------------------------
# (written in the same order the original author typed them so that any
# newly-created shared-string entries land in the same sequence)

$ws.Range("B9").Value  = "cmip6.atmos.solar.insolation_ozone.solar_ozone_impact"

$ws.Range("B11").Value = "cmip6.atmos.solar.orbital_parameters.computation_method"
$ws.Range("B12").Value = "cmip6.atmos.solar.orbital_parameters.computation_method"

$ws.Range("B15").Value = "cmip6.atmos.solar.orbital_parameters.fixed_reference_date"
$ws.Range("B16").Value = "cmip6.atmos.solar.orbital_parameters.solar_constant_transient_characteristics"
$ws.Range("B17").Value = "cmip6.atmos.solar.orbital_parameters.type"

$ws.Range("B22").Value = "cmip6.atmos.solar.solar_constant.type"
$ws.Range("B21").Value = "cmip6.atmos.solar.solar_constant.transient_characteristics"
$ws.Range("B20").Value = "cmip6.atmos.solar.solar_constant.fixed_value"
$ws.Range("B23").Value = "cmip6.atmos.solar.solar_constant.fixed_value"

$ws.Range("B25").Value = "cmip6.atmos.volcanos.volcanoes_treatment.volcanoes_implementation"

# --- Formatting touch-ups -------------------------------------------------
# B4 / B7 pick up the same body style ("s=4") used by their sibling rows.
$ws.Range("B4").Style = $ws.Range("B5").Style
$ws.Range("B7").Style = $ws.Range("B5").Style

# B15:B17 get a distinct (hyperlink-like) look: plain black Helvetica text
# on a thin grey box border.
$hlRange = $ws.Range("B15:B17")
$hlRange.Font.Color = 0
$hlRange.Borders.LineStyle = 1
$hlRange.Borders.Weight = 2
$hlRange.Borders.Color = 10855845

# Column A best-fits to the (now shorter) longest label.
$ws.Columns("A:A").AutoFit()

# Selection moves to B5, matching the author's final cursor position.
$ws.Range("B5").Select()
